$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 2: "Format:" / "v0.1.0" caption, italic Aptos 11pt ---
$ws.Range("B2").Value = "Format:"
$ws.Range("B2").Font.Italic = $true
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Name = "Aptos"
$ws.Range("B2").HorizontalAlignment = -4152
$ws.Range("B2").VerticalAlignment = -4108

$ws.Range("C2").Value = "v0.1.0"
$ws.Range("C2").Font.Italic = $true
$ws.Range("C2").Font.Size = 11
$ws.Range("C2").Font.Name = "Aptos"
$ws.Range("C2").VerticalAlignment = -4108

$ws.Rows.Item(2).RowHeight = 18.75

# --- Minor column width tweaks ---
$ws.Columns.Item(1).ColumnWidth = 4.666666666666667
$ws.Columns.Item(8).ColumnWidth = 18.5
